$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '69.562.36'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  -0.83%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.493.34'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  -0.88%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '568.96'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -1.21%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '164.81'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -0.67%  '
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.511'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -1.73%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '2.490.88'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -1.01%  '
$ws.Range('E10').Value = '  -1.55%  '
$ws.Range('E11').Value = '  -0.36%  '
$ws.Range('E12').Value = '  -0.58%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '4.91'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -0.30%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '2.943.51'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -1.28%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '69.354.58'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -0.96%  '
$ws.Range('E16').Value = '  -1.06%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '24.27'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -2.97%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '2.502.45'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -1.06%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '11.16'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -2.34%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '7.37'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -5.38%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '345.67'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -1.58%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '3.88'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -1.31%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '1.92'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -3.47%  '
$ws.Range('E24').Value = '  +0.00%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '69.61'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -1.16%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '3.91'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -2.47%  '
$ws.Range('B27').Value = 'Aptos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '8.66'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -2.17%  '
$ws.Range('B28').Value = 'WrappedeETH'
$ws.Range('C28').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '2.614.52'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -1.75%  '
$ws.Range('E29').Value = '  -0.49%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.0₃0874'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -2.71%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '7.65'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -3.01%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.19'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -5.09%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '437.19'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -5.99%  '
$ws.Range('E34').Value = '  -0.09%  '
$ws.Range('E35').Value = '  -2.23%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '154.90'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -0.61%  '
$ws.Range('E37').Value = '  -4.08%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '19.06'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -0.06%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '18.15'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -2.65%  '
$ws.Range('E40').Value = '  +0.02%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '4.59'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -3.92%  '
$ws.Range('E43').Value = '  -1.86%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '2.17'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -5.50%  '
$ws.Range('E45').Value = '  -6.95%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '138.21'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -3.20%  '
$ws.Range('E47').Value = '  -1.64%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.511'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -3.32%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.0723'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -0.82%  '
$ws.Range('E50').Value = '  -0.96%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.0921'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -0.94%  '
